$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sayfa1")

# Fill in the formulas in column E (E6:E14): total score after adding the fixed bonus in $C$6
$ws.Range("E6").Formula = "=D6+`$C`$6"
$ws.Range("E7:E14").Formula = "=D7+`$C`$6"

# Fill in student info in the I column next to the Numara/Ad Soyad/Bölüm labels
$ws.Range("I6").Value = 20215070019
$ws.Range("I7").Value = "KÜBRA ÇABUK"
$ws.Range("I8").Value = "YBS"

# Update the active selection to match the saved view
[void]$ws.Range("K11").Select()

[void]$wb.Save()
